$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the "mesh" row (row 9) to list the new canal
# network / gmsh related file pointers used by the smaller-mesh workflow.
$ws.Rows("9:12").Insert()
$ws.Range("C9:I12").Clear()

$ws.Range("A9").Value = "channel_network_nodes"
$ws.Range("B9").Value = "C:\Users\03125327\Dropbox\PhD\Computation\ForestCarbon\2022 Kalimantan customer work\qgis_derivated_data\channel_net_nodes.gpkg"

$ws.Range("A10").Value = "channel_network_lines"
$ws.Range("B10").Value = "C:/Users/03125327/Dropbox/PhD/Computation/ForestCarbon/2022 Kalimantan customer work/qgis_derivated_data/reprojected_canals_singleparts.gpkg"

$ws.Range("A11").Value = "study_area_boundary_points"
$ws.Range("B11").Value = "C:\Users\03125327\github\fc_hydro_kalimantan_2022\data\study_area_boundary_points.gpkg"

$ws.Range("A12").Value = "gmsh_geo"
$ws.Range("B12").Value = "data/mesh/mesh.geo"

# Row-insert copies formatting down from the row above, which leaves an
# unwanted style on the "Path" cells of the first three new rows (only the
# last new row, which mirrors the "Content"/"Path" pair style above it,
# should keep it). Re-assigning a no-op font property clears that inherited
# style so these cells match the unstyled "Path" cells used elsewhere.
$ws.Range("B9:B11").Font.Name = $ws.Range("B9:B11").Font.Name

# The dipwell measurements used for the initial condition now come from a
# file filtered to points far from canals.
$ws.Range("B19").Value = "initial_condition/initial_day_dipwell_coords_and_measurements_far_from_canals.csv"

# Widen column A so the longer content/path labels are not clipped.
$ws.Columns("A").ColumnWidth = 34.86

# Leave the selection where the author was last editing.
[void]$ws.Range("B13").Select()
